# Scheduled-runner refresh of Universalis market-price snapshots.
# Updates the computed Leve profit columns (H:N) for the rows whose
# backing item prices moved since the last run, across several job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2298.7058
$ws.Range("J138").Value = 2559.923
$ws.Range("L138").Value = 7679.768999999999
$ws.Range("N138").Value = -17959.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2494.1
$ws.Range("I2").Value = 1277.5714
$ws.Range("J2").Value = 5332.6665
$ws.Range("K2").Value = 1277.5714
$ws.Range("L2").Value = 5332.6665
$ws.Range("M2").Value = -1164.5714
$ws.Range("N2").Value = -5558.6665

$ws.Range("H32").Value = 5313.1665
$ws.Range("I32").Value = 5313.1665
$ws.Range("K32").Value = 5313.1665
$ws.Range("M32").Value = -5026.1665

$ws.Range("H45").Value = 3553.9092
$ws.Range("I45").Value = 1396.5
$ws.Range("J45").Value = 4033.3333
$ws.Range("K45").Value = 1396.5
$ws.Range("L45").Value = 4033.3333
$ws.Range("M45").Value = -1019.5
$ws.Range("N45").Value = -4787.3333

$ws.Range("H61").Value = 5800.5
$ws.Range("J61").Value = 16007
$ws.Range("L61").Value = 16007
$ws.Range("N61").Value = -16431

$ws.Range("H63").Value = 2631
$ws.Range("J63").Value = 2185.3333
$ws.Range("L63").Value = 2185.3333
$ws.Range("N63").Value = -3557.3333

$ws.Range("H66").Value = 2631
$ws.Range("J66").Value = 2185.3333
$ws.Range("L66").Value = 10926.6665
$ws.Range("N66").Value = -17790.6665

$ws.Range("H97").Value = 1060.174
$ws.Range("I97").Value = 938.8125
$ws.Range("J97").Value = 1337.5714
$ws.Range("K97").Value = 938.8125
$ws.Range("L97").Value = 1337.5714
$ws.Range("M97").Value = -442.8125
$ws.Range("N97").Value = -2329.5714

$ws.Range("H116").Value = 2494.1
$ws.Range("I116").Value = 1277.5714
$ws.Range("J116").Value = 5332.6665
$ws.Range("K116").Value = 1277.5714
$ws.Range("L116").Value = 5332.6665
$ws.Range("M116").Value = 1016.4286
$ws.Range("N116").Value = -9920.666499999999

$ws.Range("H134").Value = 80000
$ws.Range("J134").Value = 80000
$ws.Range("L134").Value = 80000
$ws.Range("N134").Value = -90140

$ws.Range("H136").Value = 5800.5
$ws.Range("J136").Value = 16007
$ws.Range("L136").Value = 48021
$ws.Range("N136").Value = -53121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2494.1
$ws.Range("I3").Value = 1277.5714
$ws.Range("J3").Value = 5332.6665
$ws.Range("K3").Value = 1277.5714
$ws.Range("L3").Value = 5332.6665
$ws.Range("M3").Value = -1163.5714
$ws.Range("N3").Value = -5560.6665

$ws.Range("H22").Value = 593.5
$ws.Range("I22").Value = 593.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 593.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -420.5
$ws.Range("N22").ClearContents()

$ws.Range("H99").Value = 2055.818
$ws.Range("I99").Value = 1635.1111
$ws.Range("J99").Value = 3949
$ws.Range("K99").Value = 1635.1111
$ws.Range("L99").Value = 3949
$ws.Range("M99").Value = -137.1111000000001
$ws.Range("N99").Value = -6945

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 982.5
$ws.Range("I58").Value = 982.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 982.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -779.5
$ws.Range("N58").ClearContents()

$ws.Range("H99").Value = 3356.1667
$ws.Range("I99").Value = 3356.1667
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3356.1667
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1858.1667
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 3356.1667
$ws.Range("I126").Value = 3356.1667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10068.5001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7598.500100000001
$ws.Range("N126").ClearContents()

$ws.Range("H136").Value = 982.5
$ws.Range("I136").Value = 982.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2947.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -397.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 425
$ws.Range("J17").Value = 425
$ws.Range("L17").Value = 1275
$ws.Range("N17").Value = -1613

$ws.Range("H113").Value = 618.1429000000001
$ws.Range("I113").Value = 245
$ws.Range("J113").Value = 680.3333
$ws.Range("K113").Value = 735
$ws.Range("L113").Value = 2040.9999
$ws.Range("M113").Value = 1435
$ws.Range("N113").Value = -6380.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1673666.6
$ws.Range("I3").Value = 1673666.6
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1673666.6
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1673550.6
$ws.Range("N3").ClearContents()

$ws.Range("H11").Value = 8458.333000000001

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H80").Value = 3936.3125
$ws.Range("I80").Value = 1993.5
$ws.Range("K80").Value = 1993.5
$ws.Range("M80").Value = -995.5

$ws.Range("H83").Value = 3936.3125
$ws.Range("I83").Value = 1993.5
$ws.Range("K83").Value = 9967.5
$ws.Range("M83").Value = -4975.5

$ws.Range("H126").Value = 7449.75
$ws.Range("I126").Value = 4933
$ws.Range("K126").Value = 14799
$ws.Range("M126").Value = -12329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5998.8
$ws.Range("I7").Value = 5998.8
$ws.Range("K7").Value = 5998.8
$ws.Range("M7").Value = -5886.8

$ws.Range("H22").Value = 3462.75
$ws.Range("J22").Value = 4317
$ws.Range("L22").Value = 4317
$ws.Range("N22").Value = -4907

$ws.Range("H27").Value = 3462.75
$ws.Range("J27").Value = 4317
$ws.Range("L27").Value = 4317
$ws.Range("N27").Value = -4531

$ws.Range("H40").Value = 1798.5
$ws.Range("I40").Value = 1798.5
$ws.Range("K40").Value = 1798.5
$ws.Range("M40").Value = -1662.5

$ws.Range("H46").Value = 3950
$ws.Range("I46").Value = 3875
$ws.Range("J46").Value = 4100
$ws.Range("K46").Value = 3875
$ws.Range("L46").Value = 4100
$ws.Range("M46").Value = -3687
$ws.Range("N46").Value = -4476

$ws.Range("H68").Value = 7000
$ws.Range("I68").Value = 7000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 7000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -6251
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 7000
$ws.Range("I71").Value = 7000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 35000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -31256
$ws.Range("N71").ClearContents()

$ws.Range("H126").Value = 5998.8
$ws.Range("I126").Value = 5998.8
$ws.Range("K126").Value = 17996.4
$ws.Range("M126").Value = -15526.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 938.8570999999999
$ws.Range("I136").Value = 857.2308
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 2571.6924
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -21.69239999999991
$ws.Range("N136").Value = -11100
